# Fruta / hortaliza, semanal
#
# A new daily price record was added to the "Fruta, Macroferia Regional de
# Talca - Ciruela" table. In the canonical OOXML this shows up as a brand
# new row 10 (the former rows 10..101 all shift down by one, to 11..102,
# and the sheet's used range grows from A1:T101 to A1:T102).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the current row 10 — this pushes the
# existing rows 10-101 down to 11-102 (matching the diff exactly).
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new record's data.
$ws.Range("A10").Value = 5
$ws.Range("B10").Value = "Macroferia Regional de Talca"
$ws.Range("C10").Value = "Maule"
$ws.Range("D10").Value = 44635
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100103
$ws.Range("H10").Value = "Frutos de hueso (carozo)"
$ws.Range("I10").Value = 100103002
$ws.Range("J10").Value = "Ciruela"
$ws.Range("K10").Value = "Angeleno"
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 6000
$ws.Range("O10").Value = 6000
$ws.Range("P10").Value = 6000
$ws.Range("Q10").Value = "$/bandeja 18 kilos granel"
$ws.Range("R10").Value = "Provincia de Curicó"
$ws.Range("S10").Value = 333
$ws.Range("T10").Value = 18
